# Applies the scheduled market-data refresh to the Exodus_Profits workbook.
# For each Leve row, refreshes currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and
# the derived LeveProfit(NQ/HQ) columns (H, I, J, K, L, M, N) to the latest
# market-board snapshot. A couple of rows gain/lose their HQ-profit cell (N)
# entirely when an HQ price did/did not come back from the board this run.

$wb = $excel.ActiveWorkbook
$totalUpdates = 0

# ---- Sheet: ALC (56 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 312.33334
$ws.Range("I2").Value = 648.75
$ws.Range("J2").Value = 43.2
$ws.Range("K2").Value = 648.75
$ws.Range("L2").Value = 43.2
$ws.Range("M2").Value = -535.75
$ws.Range("N2").Value = -269.2
$ws.Range("H4").Value = 657.7143
$ws.Range("I4").Value = 657.7143
$ws.Range("K4").Value = 657.7143
$ws.Range("M4").Value = -543.7143
$ws.Range("H28").Value = 3987.5557
$ws.Range("I28").Value = 2383.8
$ws.Range("J28").Value = 5992.25
$ws.Range("K28").Value = 2383.8
$ws.Range("L28").Value = 5992.25
$ws.Range("M28").Value = -1898.8
$ws.Range("N28").Value = -6962.25
$ws.Range("H74").Value = 4480.2
$ws.Range("I74").Value = 4444.6665
$ws.Range("K74").Value = 4444.6665
$ws.Range("M74").Value = -3508.6665
$ws.Range("H77").Value = 4480.2
$ws.Range("I77").Value = 4444.6665
$ws.Range("K77").Value = 22223.3325
$ws.Range("M77").Value = -17543.3325
$ws.Range("H110").Value = 67898.336
$ws.Range("J110").Value = 67898.336
$ws.Range("L110").Value = 67898.336
$ws.Range("N110").Value = -76078.336
$ws.Range("H112").Value = 1169.2858
$ws.Range("I112").Value = 695.1667
$ws.Range("J112").Value = 1298.591
$ws.Range("K112").Value = 2085.5001
$ws.Range("L112").Value = 3895.773
$ws.Range("M112").Value = -977.5001000000002
$ws.Range("N112").Value = -6111.772999999999
$ws.Range("H121").Value = 2000
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -9494
$ws.Range("H134").Value = 61014.168
$ws.Range("J134").Value = 61014.168
$ws.Range("L134").Value = 61014.168
$ws.Range("N134").Value = -71154.16800000001
$ws.Range("H136").Value = 68096.10000000001
$ws.Range("J136").Value = 68096.10000000001
$ws.Range("L136").Value = 68096.10000000001
$ws.Range("N136").Value = -78296.10000000001
$ws.Range("H138").Value = 2010
$ws.Range("I138").Value = 1744.683
$ws.Range("J138").Value = 2735.2
$ws.Range("K138").Value = 5234.049
$ws.Range("L138").Value = 8205.599999999999
$ws.Range("M138").Value = -94.04899999999998
$ws.Range("N138").Value = -18485.6
$totalUpdates += 56

# ---- Sheet: ARM (40 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 227.375
$ws.Range("I4").Value = 234.14285
$ws.Range("J4").Value = 180
$ws.Range("K4").Value = 234.14285
$ws.Range("L4").Value = 180
$ws.Range("M4").Value = -118.14285
$ws.Range("N4").Value = -412
$ws.Range("H6").Value = 535000
$ws.Range("J6").Value = 80000
$ws.Range("L6").Value = 80000
$ws.Range("N6").Value = -80346
$ws.Range("H7").Value = 27530.334
$ws.Range("J7").Value = 27530.334
$ws.Range("L7").Value = 27530.334
$ws.Range("N7").Value = -27758.334
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H25").Value = 135
$ws.Range("H61").Value = 1786.878
$ws.Range("I61").Value = 1570.7838
$ws.Range("K61").Value = 1570.7838
$ws.Range("M61").Value = -1358.7838
$ws.Range("H107").Value = 39204.57
$ws.Range("J107").Value = 39204.57
$ws.Range("L107").Value = 39204.57
$ws.Range("N107").Value = -46884.57
$ws.Range("H118").Value = 47331.332
$ws.Range("J118").Value = 47331.332
$ws.Range("L118").Value = 47331.332
$ws.Range("N118").Value = -50645.332
$ws.Range("H136").Value = 1786.878
$ws.Range("I136").Value = 1570.7838
$ws.Range("K136").Value = 4712.3514
$ws.Range("M136").Value = -2162.3514
$totalUpdates += 40

# ---- Sheet: BSM (37 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47229.184
$ws.Range("I20").Value = 68221.664
$ws.Range("K20").Value = 68221.664
$ws.Range("M20").Value = -67974.664
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H109").Value = 96658.164
$ws.Range("J109").Value = 96658.164
$ws.Range("L109").Value = 96658.164
$ws.Range("N109").Value = -99432.164
$ws.Range("H115").Value = 84397
$ws.Range("H119").Value = 91658.164
$ws.Range("J119").Value = 91658.164
$ws.Range("L119").Value = 91658.164
$ws.Range("N119").Value = -101334.164
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H127").Value = 58553.145
$ws.Range("J127").Value = 58553.145
$ws.Range("L127").Value = 58553.145
$ws.Range("N127").Value = -68473.14499999999
$ws.Range("H132").Value = 89309.5
$ws.Range("J132").Value = 89309.5
$ws.Range("L132").Value = 89309.5
$ws.Range("N132").Value = -99429.5
$ws.Range("H135").Value = 58612.8
$ws.Range("J135").Value = 58612.8
$ws.Range("L135").Value = 58612.8
$ws.Range("N135").Value = -68752.8
$ws.Range("H138").Value = 80831.375
$ws.Range("J138").Value = 80831.375
$ws.Range("L138").Value = 80831.375
$ws.Range("N138").Value = -91111.375
$totalUpdates += 37

# ---- Sheet: CRP (20 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 48424.5
$ws.Range("J18").Value = 48424.5
$ws.Range("L18").Value = 48424.5
$ws.Range("N18").Value = -48884.5
$ws.Range("H31").Value = 2293.6924
$ws.Range("I31").Value = 1828.4286
$ws.Range("K31").Value = 1828.4286
$ws.Range("M31").Value = -1533.4286
$ws.Range("H34").Value = 2293.6924
$ws.Range("I34").Value = 1828.4286
$ws.Range("K34").Value = 1828.4286
$ws.Range("M34").Value = -1626.4286
$ws.Range("H114").Value = 44984.125
$ws.Range("J114").Value = 44984.125
$ws.Range("L114").Value = 44984.125
$ws.Range("N114").Value = -53662.125
$ws.Range("H118").Value = 59842.223
$ws.Range("J118").Value = 59842.223
$ws.Range("L118").Value = 59842.223
$ws.Range("N118").Value = -63156.223
$totalUpdates += 20

# ---- Sheet: CUL (7 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 6337.6
$ws.Range("I132").Value = 4090.75
$ws.Range("J132").Value = 7835.5
$ws.Range("K132").Value = 36816.75
$ws.Range("L132").Value = 70519.5
$ws.Range("M132").Value = -34286.75
$ws.Range("N132").Value = -75579.5
$totalUpdates += 7

# ---- Sheet: GSM (34 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 486.35715
$ws.Range("I2").Value = 299.42856
$ws.Range("J2").Value = 673.2857
$ws.Range("K2").Value = 299.42856
$ws.Range("L2").Value = 673.2857
$ws.Range("M2").Value = -186.42856
$ws.Range("N2").Value = -899.2857
$ws.Range("H80").Value = 3620.4285
$ws.Range("J80").Value = 3501.3333
$ws.Range("L80").Value = 3501.3333
$ws.Range("N80").Value = -5497.3333
$ws.Range("H83").Value = 3620.4285
$ws.Range("J83").Value = 3501.3333
$ws.Range("L83").Value = 17506.6665
$ws.Range("N83").Value = -27490.6665
$ws.Range("H110").Value = 74911.91
$ws.Range("J110").Value = 74911.91
$ws.Range("L110").Value = 74911.91
$ws.Range("N110").Value = -83091.91
$ws.Range("H113").Value = 1588957.2
$ws.Range("I113").Value = 1677.875
$ws.Range("J113").Value = 6668251
$ws.Range("K113").Value = 1677.875
$ws.Range("L113").Value = 6668251
$ws.Range("M113").Value = 492.125
$ws.Range("N113").Value = -6672591
$ws.Range("H135").Value = 50450.41
$ws.Range("J135").Value = 50450.41
$ws.Range("L135").Value = 50450.41
$ws.Range("N135").Value = -60590.41
$ws.Range("H140").Value = 98496.664
$ws.Range("J140").Value = 98496.664
$ws.Range("L140").Value = 98496.664
$ws.Range("N140").Value = -108856.664
$totalUpdates += 34

# ---- Sheet: LTW (20 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1359.2
$ws.Range("I93").Value = 1359.2
$ws.Range("K93").Value = 1359.2
$ws.Range("M93").Value = -111.2
$ws.Range("H100").Value = 10148.708
$ws.Range("I100").Value = 10488.353
$ws.Range("K100").Value = 10488.353
$ws.Range("M100").Value = -9947.352999999999
$ws.Range("H118").Value = 53636.727
$ws.Range("J118").Value = 53636.727
$ws.Range("L118").Value = 53636.727
$ws.Range("N118").Value = -56950.727
$ws.Range("H123").Value = 79108.89
$ws.Range("J123").Value = 79108.89
$ws.Range("L123").Value = 79108.89
$ws.Range("N123").Value = -88908.89
$ws.Range("H136").Value = 2353.45
$ws.Range("J136").Value = 2678.875
$ws.Range("L136").Value = 8036.625
$ws.Range("N136").Value = -13136.625
$totalUpdates += 20

# ---- Sheet: WVR (15 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1762.7561
$ws.Range("I107").Value = 1206.9584
$ws.Range("J107").Value = 2547.4119
$ws.Range("K107").Value = 3620.8752
$ws.Range("L107").Value = 7642.2357
$ws.Range("M107").Value = -1700.8752
$ws.Range("N107").Value = -11482.2357
$ws.Range("H127").Value = 92078
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920
$ws.Range("H136").Value = 1470.1666
$ws.Range("I136").Value = 886.3158
$ws.Range("K136").Value = 2658.9474
$ws.Range("M136").Value = -108.9474
$totalUpdates += 15

Write-Output "Updated $totalUpdates cells across $($wb.Worksheets.Count) sheets"
